# Update countries & provincias Spain
# Applies the 3-Apr-2020 18:50 -> 19:20 data refresh to the "paises" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / timestamp text
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 19:20"

# ---------------------------------------------------------------------------
# Simple in-place numeric refreshes (row/country unchanged)
# Columns: A Pais | B Casos totales | C Nuevos casos | D Casos activos
#          E Recuperados | F Casos criticos | G Muertes hoy | H Muertes
# ---------------------------------------------------------------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 265506
$ws.Range("C4").Value = 20629
$ws.Range("D4").Value = 11983
$ws.Range("E4").Value = 246737
$ws.Range("G4").Value = 716
$ws.Range("H4").Value = 6786

# Row 7 - Alemania
$ws.Range("B7").Value = 89838
$ws.Range("C7").Value = 5044
$ws.Range("E7").Value = 64033
$ws.Range("G7").Value = 123
$ws.Range("H7").Value = 1230

# Row 13 - Suiza
$ws.Range("B13").Value = 19606
$ws.Range("C13").Value = 779
$ws.Range("E13").Value = 14169
$ws.Range("G13").Value = 55
$ws.Range("H13").Value = 591

# Row 17 - Austria
$ws.Range("B17").Value = 11489
$ws.Range("C17").Value = 360
$ws.Range("E17").Value = 9299

# Row 33 - Chile
$ws.Range("E33").Value = 2767
$ws.Range("G33").Value = 18
$ws.Range("H33").Value = 133

# Row 35 - Malasia
$ws.Range("E35").Value = 2465
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 40

# Row 72 - Bosnia y Herzegovina
$ws.Range("B72").Value = 575
$ws.Range("C72").Value = 42
$ws.Range("E72").Value = 531

# Row 93
$ws.Range("B93").Value = 281
$ws.Range("C93").Value = 8
$ws.Range("E93").Value = 265

# Row 109
$ws.Range("B109").Value = 159
$ws.Range("C109").Value = 8
$ws.Range("E109").Value = 131

# Row 125
$ws.Range("D125").Value = 6
$ws.Range("E125").Value = 83

# ---------------------------------------------------------------------------
# Armenia / Marruecos swap (rows 65-66): Marruecos' count overtakes Armenia's
# ---------------------------------------------------------------------------
$ws.Range("A65").Value = "Marruecos"
$ws.Range("B65").Value = 761
$ws.Range("C65").Value = 53
$ws.Range("D65").Value = 56
$ws.Range("E65").Value = 658
$ws.Range("F65").Value = 1
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 47

$ws.Range("A66").Value = "Armenia"
$ws.Range("B66").Value = 736
$ws.Range("C66").Value = 73
$ws.Range("D66").Value = 43
$ws.Range("E66").Value = 686
$ws.Range("F66").Value = 30
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 7

# ---------------------------------------------------------------------------
# Guinea / Madagascar / Aruba / Banglades / Monaco / Guayana Francesa
# reshuffle (rows 128-133): Guinea's update moves it above Madagascar,
# Aruba's update moves it above Banglades, pushing the rest down one slot.
# ---------------------------------------------------------------------------
$ws.Range("A128").Value = "Guinea"
$ws.Range("B128").Value = 73
$ws.Range("C128").Value = 21
$ws.Range("D128").Value = 2
$ws.Range("E128").Value = 71
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 0

$ws.Range("A129").Value = "Madagascar"
$ws.Range("B129").Value = 65
$ws.Range("C129").Value = 6
$ws.Range("D129").Value = 0
$ws.Range("E129").Value = 65
$ws.Range("F129").Value = 6
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 0

$ws.Range("A130").Value = "Aruba"
$ws.Range("B130").Value = 62
$ws.Range("C130").Value = 2
$ws.Range("D130").Value = 1
$ws.Range("E130").Value = 61
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 0

$ws.Range("A131").Value = "Banglades"
$ws.Range("B131").Value = 61
$ws.Range("C131").Value = 5
$ws.Range("D131").Value = 26
$ws.Range("E131").Value = 29
$ws.Range("F131").Value = 1
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 6

$ws.Range("A132").Value = "Monaco"
$ws.Range("B132").Value = 60
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 2
$ws.Range("E132").Value = 57
$ws.Range("F132").Value = 2
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 1

$ws.Range("A133").Value = "Guayana Francesa"
$ws.Range("B133").Value = 57
$ws.Range("C133").Value = 6
$ws.Range("D133").Value = 22
$ws.Range("E133").Value = 35
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0
